$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.831.77"
$ws.Range("D3").Value = "3.404.30"
$ws.Range("E3").Value = "  +3.21%  "
$ws.Range("E4").Value = "  +0.04%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "578.50"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("E6").Value = "  +8.17%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.405.00"
$ws.Range("E8").Value = "  +3.24%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.478"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "7.69"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("E11").Value = "  +5.94%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.395"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +5.50%  "
$ws.Range("D13").Value = "3.981.95"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("E15").Value = "  +6.49%  "
$ws.Range("D16").Value = "3.440.48"
$ws.Range("E16").Value = "  +3.88%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "25.57"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +4.91%  "
$ws.Range("D18").Value = "61.955.87"
$ws.Range("E18").Value = "  +1.25%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "14.03"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +5.32%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "5.91"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +4.69%  "
$ws.Range("E21").Value = "  +6.09%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "389.97"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +10.18%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "0.572"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("D24").Value = "3.534.67"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("E25").Value = "  +0.14%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "71.20"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +2.97%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "0.0000125"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +16.69%  "
$ws.Range("E28").Value = "  +15.53%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "7.83"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +8.98%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "0.994"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "8.27"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +5.99%  "
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "3.437.20"
$ws.Range("E35").Value = "  +3.22%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "23.63"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("E37").Value = "  +6.02%  "
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("E39").Value = "  +5.68%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "161.76"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.0794"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +5.14%  "
$ws.Range("E43").Value = "  +0.06%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "1.22"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +9.92%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "4.46"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +2.12%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.772"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("E47").Value = "  +1.07%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "23.77"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +6.92%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "7.03"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +5.26%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "22.99"
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +8.78%  "
$ws.Range("D51").Value = "2.344.98"
$ws.Range("E51").Value = "  +8.77%  "

Write-Output "Applied 80 changes"
